$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-18 down to 17-19
$ws.Rows("16:16").Insert()

# Fill in the new row 16 with data (matching the format of its neighboring rows)
$ws.Cells.Item(16, 1).Value = 4
$ws.Cells.Item(16, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(16, 3).Value = "Los Lagos"
$ws.Cells.Item(16, 4).Value = 44875
$ws.Cells.Item(16, 5).Value = 10
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100101
$ws.Cells.Item(16, 8).Value = "Berries"
$ws.Cells.Item(16, 9).Value = 100101001
$ws.Cells.Item(16, 10).Value = "Arándano (blue)"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 400
$ws.Cells.Item(16, 14).Value = 7500
$ws.Cells.Item(16, 15).Value = 8000
$ws.Cells.Item(16, 16).Value = 7750
$ws.Cells.Item(16, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(16, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(16, 19).Value = 5167
$ws.Cells.Item(16, 20).Value = 1.5
